$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh Price (column D) and Volume 1h change (column E) for each crypto row
# with the latest scraped values. D-column values are forced to text via a
# quote-prefix so numeric-looking strings (e.g. "1.003", "45.20") are not
# reinterpreted as numbers, then the cell style is restored to "Normal".
$ws.Range("D2").Value = "'27.740.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "'1.850.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'313.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.4330"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.60%  "
$ws.Range("D8").Value = "'0.3658"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("D9").Value = "'45.20"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("D10").Value = "'0.07334"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.35%  "
$ws.Range("D11").Value = "'0.8801"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.28%  "
$ws.Range("D12").Value = "'20.75"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "'1.834.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").Value = "'5.344"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("D15").Value = "'6.525"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("D16").Value = "'0.06929"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.69%  "
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "'80.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.56%  "
$ws.Range("D19").Value = "'0.000009046"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.91%  "
$ws.Range("D20").Value = "'1.004"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").Value = "'15.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").Value = "'28.007.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.27%  "
$ws.Range("D23").Value = "'4.978"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("D24").Value = "'10.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.54%  "
$ws.Range("D25").Value = "'2.167.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.75%  "
$ws.Range("D26").Value = "'1.991"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.66%  "
$ws.Range("D27").Value = "'155.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("D28").Value = "'18.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.77%  "
$ws.Range("D29").Value = "'120.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.86%  "
$ws.Range("D30").Value = "'5.266"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").Value = "'1.856"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.26%  "
$ws.Range("D32").Value = "'0.08928"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").Value = "'0.7608"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.26%  "
$ws.Range("D34").Value = "'4.554"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("D35").Value = "'2.949"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.55%  "
$ws.Range("D36").Value = "'1.124"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.60%  "
$ws.Range("D37").Value = "'1.109"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.23%  "
$ws.Range("D38").Value = "'0.05414"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("D39").Value = "'0.01938"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.77%  "
$ws.Range("D40").Value = "'2.842"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.47%  "
$ws.Range("D41").Value = "'0.5096"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("D42").Value = "'0.1662"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.34%  "
$ws.Range("D43").Value = "'6.677"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.18%  "
$ws.Range("D44").Value = "'8.341"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.91%  "
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("D46").Value = "'0.06551"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.32%  "
$ws.Range("D47").Value = "'0.4674"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.03%  "
$ws.Range("D48").Value = "'104.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("D49").Value = "'1.000"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").Value = "'1.621"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.07%  "
$ws.Range("D51").Value = "'64.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.11%  "
